$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 584 (shifts old rows 584-635 down to 585-636)
$ws.Rows(584).Insert()

# Populate the newly inserted row 584 with values (same categorical data as the
# row that used to be 584, now shifted to 585, but with new measurement values)
$ws.Cells.Item(584, 1).Value = 10
$ws.Cells.Item(584, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(584, 3).Value = "La Araucanía"
$ws.Cells.Item(584, 4).Value = 45166
$ws.Cells.Item(584, 5).Value = 9
$ws.Cells.Item(584, 6).Value = "Fruta"
$ws.Cells.Item(584, 7).Value = 100108
$ws.Cells.Item(584, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(584, 9).Value = 100108002
$ws.Cells.Item(584, 10).Value = "Mango"
$ws.Cells.Item(584, 11).Value = "Sin especificar"
$ws.Cells.Item(584, 12).Value = "Primera"
$ws.Cells.Item(584, 13).Value = 505
$ws.Cells.Item(584, 14).Value = 10000
$ws.Cells.Item(584, 15).Value = 11000
$ws.Cells.Item(584, 16).Value = 10446
$ws.Cells.Item(584, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(584, 18).Value = "Brasil"
$ws.Cells.Item(584, 19).Value = 2612
$ws.Cells.Item(584, 20).Value = 4
